$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source row 16 ("SITIO GRANDE 2") is removed entirely; Excel shifts
# every row below it up by one (EntireRow delete).
$ws.Rows(16).Delete()

# Column B (terminación date, stored as date serials) gets a custom
# date display format, which also widens the column via autofit.
$ws.Columns(2).NumberFormat = "dd\-mm\-yy;@"

# Selection returns to the top of the sheet, landing on the row that is
# now 16 (used to be 17).
$ws.Range("A16").Select()
